$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.638.29"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3
$ws.Range("D3").Value = "2.342.71"
$ws.Range("E3").Value = "  -0.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'305.37"
$ws.Range("E5").Value = "  -1.53%  "

# Row 6
$ws.Range("D6").Value = "'102.14"
$ws.Range("E6").Value = "  -1.96%  "

# Row 7
$ws.Range("E7").Value = "  -2.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  +0.52%  "

# Row 10
$ws.Range("D10").Value = "'35.45"
$ws.Range("E10").Value = "  -1.71%  "

# Row 11
$ws.Range("D11").Value = "'51.66"
$ws.Range("E11").Value = "  -2.86%  "

# Row 12
$ws.Range("E12").Value = "  -1.57%  "

# Row 13
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("D14").Value = "'6.84"
$ws.Range("E14").Value = "  -2.23%  "

# Row 15
$ws.Range("D15").Value = "2.702.86"

# Row 16
$ws.Range("D16").Value = "'15.61"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("D17").Value = "2.346.73"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18
$ws.Range("D18").Value = "'0.810"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").Value = "43.529.85"
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").Value = "'11.84"
$ws.Range("E20").Value = "  -1.19%  "

# Row 21
$ws.Range("E21").Value = "  -1.56%  "

# Row 22
$ws.Range("D22").Value = "'6.14"
$ws.Range("E22").Value = "  -1.84%  "

# Row 23
$ws.Range("D23").Value = "'68.19"
$ws.Range("E23").Value = "  -0.24%  "

# Row 24
$ws.Range("D24").Value = "'238.86"
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
$ws.Range("E25").Value = "  -3.09%  "

# Row 26
$ws.Range("D26").Value = "'2.55"
$ws.Range("E26").Value = "  -3.39%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("D28").Value = "'25.10"
$ws.Range("E28").Value = "  -2.27%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'34.76"
$ws.Range("E30").Value = "  -5.86%  "

# Row 31
$ws.Range("D31").Value = "'9.29"
$ws.Range("E31").Value = "  -3.10%  "

# Row 32
$ws.Range("D32").Value = "'165.62"
$ws.Range("E32").Value = "  +2.13%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.17%  "

# Row 34
$ws.Range("D34").Value = "'5.08"
$ws.Range("E34").Value = "  -3.61%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'4.54"
$ws.Range("E35").Value = "  -4.10%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = "  -4.98%  "

# Row 37
$ws.Range("D37").Value = "'17.08"
$ws.Range("E37").Value = "  -6.44%  "

# Row 38
$ws.Range("D38").Value = "'0.0709"
$ws.Range("E38").Value = "  -4.37%  "

# Row 39
$ws.Range("E39").Value = "  -7.02%  "

# Row 40
$ws.Range("D40").Value = "'1.84"
$ws.Range("E40").Value = "  -5.10%  "

# Row 41
$ws.Range("E41").Value = "  -2.48%  "

# Row 42
$ws.Range("E42").Value = "  -2.45%  "

# Row 43
$ws.Range("D43").Value = "'2.45"
$ws.Range("E43").Value = "  -1.04%  "

# Row 44
$ws.Range("D44").Value = "1.992.37"
$ws.Range("E44").Value = "  -0.42%  "

# Row 45
$ws.Range("E45").Value = "  -1.37%  "

# Row 46
$ws.Range("D46").Value = "'18.68"
$ws.Range("E46").Value = "  -8.21%  "

# Row 47
$ws.Range("E47").Value = "  -5.99%  "

# Row 48
$ws.Range("D48").Value = "'9.96"
$ws.Range("E48").Value = "  -4.30%  "

# Row 49
$ws.Range("E49").Value = "  -2.79%  "

# Row 50
$ws.Range("D50").Value = "'4.93"
$ws.Range("E50").Value = "  +4.86%  "

# Row 51
$ws.Range("D51").Value = "2.567.46"
$ws.Range("E51").Value = "  -1.39%  "
